# Apply updated cryptocurrency price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to be stored as text (matching the original inline-string
    # cell type) even when the text looks like a number (e.g. "0.998", "2.20").
    $range.NumberFormat = "@"
    $range.Value = $text
    # Reset the style back to Normal so no stray number-format style is left
    # on the cell (keeps cell formatting identical to the original workbook).
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "47.796.74"
$ws.Range("E2").Value = "  -0.84%  "
Set-TextValue $ws.Range("D3") "2.478.96"
$ws.Range("E3").Value = "  -1.71%  "
Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue $ws.Range("D5") "316.66"
$ws.Range("E5").Value = "  -2.01%  "
Set-TextValue $ws.Range("D6") "104.96"
$ws.Range("E6").Value = "  -4.68%  "
Set-TextValue $ws.Range("D7") "0.518"
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("E8").Value = "  -0.04%  "
Set-TextValue $ws.Range("D9") "0.534"
$ws.Range("E9").Value = "  -3.85%  "
Set-TextValue $ws.Range("D10") "38.88"
$ws.Range("E10").Value = "  -4.33%  "
Set-TextValue $ws.Range("D11") "20.17"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("E14").Value = "  -3.37%  "
Set-TextValue $ws.Range("D15") "2.864.93"
$ws.Range("E15").Value = "  -1.81%  "
Set-TextValue $ws.Range("D16") "2.489.95"
$ws.Range("E16").Value = "  -1.19%  "
Set-TextValue $ws.Range("D17") "0.823"
$ws.Range("E17").Value = "  -3.86%  "
Set-TextValue $ws.Range("D18") "47.696.83"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").Value = "  +9.45%  "
Set-TextValue $ws.Range("D20") "12.66"
$ws.Range("E20").Value = "  -4.28%  "
Set-TextValue $ws.Range("D21") "6.51"
$ws.Range("E21").Value = "  -1.75%  "
Set-TextValue $ws.Range("D22") "0.0₃0926"
$ws.Range("E22").Value = "  -2.52%  "
Set-TextValue $ws.Range("D23") "272.38"
$ws.Range("E23").Value = "  +3.17%  "
Set-TextValue $ws.Range("D24") "70.53"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("E25").Value = "  -3.76%  "
$ws.Range("E26").Value = "  +0.22%  "
Set-TextValue $ws.Range("D27") "25.57"
$ws.Range("E27").Value = "  -2.39%  "
Set-TextValue $ws.Range("D28") "2.20"
$ws.Range("E28").Value = "  -0.31%  "
Set-TextValue $ws.Range("D29") "9.59"
$ws.Range("E29").Value = "  -5.64%  "
Set-TextValue $ws.Range("D30") "0.136"
$ws.Range("E30").Value = "  -5.82%  "
Set-TextValue $ws.Range("D31") "34.50"
$ws.Range("E31").Value = "  -5.56%  "
Set-TextValue $ws.Range("D32") "49.12"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("E33").Value = "  -0.23%  "
Set-TextValue $ws.Range("D34") "18.97"
$ws.Range("E34").Value = "  -5.18%  "
Set-TextValue $ws.Range("D35") "5.24"
$ws.Range("E35").Value = "  -3.08%  "
Set-TextValue $ws.Range("D36") "0.0769"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("E37").Value = "  -3.05%  "
Set-TextValue $ws.Range("D38") "4.52"
$ws.Range("E38").Value = "  -4.28%  "
$ws.Range("E39").Value = "  -5.24%  "
Set-TextValue $ws.Range("D40") "122.64"
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("E41").Value = "  -2.23%  "
Set-TextValue $ws.Range("D42") "2.20"
$ws.Range("E42").Value = "  +0.30%  "
Set-TextValue $ws.Range("D43") "21.78"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  -0.14%  "
Set-TextValue $ws.Range("D45") "1.991.78"
$ws.Range("E45").Value = "  -1.46%  "
Set-TextValue $ws.Range("D46") "3.15"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("E47").Value = "  -1.73%  "
Set-TextValue $ws.Range("D48") "1.98"
$ws.Range("E48").Value = "  -2.85%  "
Set-TextValue $ws.Range("D49") "8.87"
$ws.Range("E49").Value = "  -2.84%  "
Set-TextValue $ws.Range("D50") "5.11"
$ws.Range("E50").Value = "  -3.00%  "
Set-TextValue $ws.Range("D51") "78.29"
$ws.Range("E51").Value = "  -0.98%  "
